$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..529).
# All of these were bumped from 45186 (2023-09-17) to 45188 (2023-09-19).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 529 }

$ws.Range("C2:C$lastRow").Value = 45188
